# Update odds data for the matches in rows 2, 5 and 6 of the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 2.7
$ws.Range("I2").Value = 3.75
$ws.Range("L2").Value = 4.35
$ws.Range("N2").Value = 5.1
$ws.Range("Q2").Value = 2.55
$ws.Range("R2").Value = 1.45
$ws.Range("S2").Value = 4.55
$ws.Range("Y2").Value = 5.6
$ws.Range("Z2").Value = 9.5
$ws.Range("AA2").Value = 9.5
$ws.Range("AB2").Value = 23
$ws.Range("AC2").Value = 23
$ws.Range("AE2").Value = 5.1
$ws.Range("AF2").Value = 5.4
$ws.Range("AH2").Value = 100
$ws.Range("AJ2").Value = 8.25
$ws.Range("AK2").Value = 19

# --- Row 5 ---
$ws.Range("G5").Value = 1.36
$ws.Range("H5").Value = 5.5
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 1.73
$ws.Range("K5").Value = 2.75
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 23
$ws.Range("N5").Value = 1.03
$ws.Range("Q5").Value = 1.36
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 1.91
$ws.Range("T5").Value = 1.8
$ws.Range("AB5").Value = 10
$ws.Range("AE5").Value = 23
$ws.Range("AJ5").Value = 26
$ws.Range("AL5").Value = 21
$ws.Range("AM5").Value = 81

# --- Row 6 ---
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 8
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.1
$ws.Range("P6").Value = 6.5
$ws.Range("Q6").Value = 1.33
$ws.Range("R6").Value = 3.25
$ws.Range("S6").Value = 1.83
$ws.Range("T6").Value = 1.83
$ws.Range("U6").Value = 1.2
$ws.Range("V6").Value = 4.33
$ws.Range("W6").Value = 1.83
$ws.Range("X6").Value = 1.83
$ws.Range("Y6").Value = 12
$ws.Range("Z6").Value = 8
$ws.Range("AB6").Value = 8
$ws.Range("AD6").Value = 23
$ws.Range("AE6").Value = 23
$ws.Range("AF6").Value = 15
$ws.Range("AH6").Value = 51
$ws.Range("AI6").Value = 700
$ws.Range("AJ6").Value = 34
$ws.Range("AN6").Value = 51
